# Insert a new data row at row 30 (shifting existing rows 30-60 down to 31-61)
# and populate it with a new weekly price record for Membrillo (Primera quality).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(30).Insert()

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C30").Value = "Metropolitana"
$ws.Range("D30").Value = 44629
$ws.Range("E30").Value = 13
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100104
$ws.Range("H30").Value = "Frutos de pepita"
$ws.Range("I30").Value = 100104003
$ws.Range("J30").Value = "Membrillo"
$ws.Range("K30").Value = "Champion"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 14
$ws.Range("N30").Value = 310000
$ws.Range("O30").Value = 320000
$ws.Range("P30").Value = 315714
$ws.Range("Q30").Value = "$/bins (450 kilos)"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 702
$ws.Range("T30").Value = 450
